$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Cells whose new text value would otherwise be auto-detected as a number by Excel
# (the source data keeps these as plain text / inline strings), so force Text format first.
$textCells = @("D4","D5","D6","D7","D8","D9","D10","D11","D13","D16","D18","D23","D25","D26","D28","D31","D32","D33","D34","D35","D37","D38","D39","D40","D42","D43","D44","D45","D46","D47","D48","D50")
foreach ($cellRef in $textCells) {
    $ws.Range($cellRef).NumberFormat = "@"
}

$ws.Range("D2").Value = "28.746.51"
$ws.Range("E2").Value = "  +2.53%  "
$ws.Range("D3").Value = "1.873.71"
$ws.Range("D4").Value = "1.005"
$ws.Range("D5").Value = "324.35"
$ws.Range("E5").Value = "  +0.03%  "
$ws.Range("D6").Value = "1.003"
$ws.Range("E6").Value = "  +0.24%  "
$ws.Range("D7").Value = "0.4605"
$ws.Range("E7").Value = "  -0.86%  "
$ws.Range("D8").Value = "0.3872"
$ws.Range("E8").Value = "  +0.11%  "
$ws.Range("D9").Value = "0.07859"
$ws.Range("E9").Value = "  +0.25%  "
$ws.Range("D10").Value = "0.9883"
$ws.Range("E10").Value = "  +2.93%  "
$ws.Range("D11").Value = "21.80"
$ws.Range("E11").Value = "  -0.48%  "
$ws.Range("D12").Value = "1.857.48"
$ws.Range("E12").Value = "  +2.28%  "
$ws.Range("D13").Value = "6.990"
$ws.Range("E13").Value = "  +1.43%  "
$ws.Range("E14").Value = "  +0.44%  "
$ws.Range("E15").Value = "  +1.86%  "
$ws.Range("D16").Value = "88.45"
$ws.Range("E16").Value = "  +0.30%  "
$ws.Range("E17").Value = "  +0.39%  "
$ws.Range("D18").Value = "0.00001004"
$ws.Range("E18").Value = "  +1.22%  "
$ws.Range("E19").Value = "  +0.82%  "
$ws.Range("D21").Value = "28.744.66"
$ws.Range("E21").Value = "  +2.48%  "
$ws.Range("E22").Value = "  -0.25%  "
$ws.Range("D23").Value = "11.05"
$ws.Range("E23").Value = "  +0.78%  "
$ws.Range("B24").Value = "WrappedliquidstakedEther2.0"
$ws.Range("C24").Value = "https://coinranking.com/coin/CiixT63n3+wrappedliquidstakedether20-wsteth"
$ws.Range("D24").Value = "2.096.86"
$ws.Range("E24").Value = "  +2.86%  "
$ws.Range("B25").Value = "Toncoin"
$ws.Range("C25").Value = "https://coinranking.com/coin/67YlI0K1b+toncoin-ton"
$ws.Range("D25").Value = "2.099"
$ws.Range("E25").Value = "  +0.43%  "
$ws.Range("D26").Value = "152.92"
$ws.Range("E26").Value = "  -1.27%  "
$ws.Range("E27").Value = "  +0.51%  "
$ws.Range("D28").Value = "5.864"
$ws.Range("E28").Value = "  +3.72%  "
$ws.Range("E29").Value = "  +0.96%  "
$ws.Range("E30").Value = "  +0.51%  "
$ws.Range("D31").Value = "0.09319"
$ws.Range("E31").Value = "  +0.78%  "
$ws.Range("D32").Value = "0.9202"
$ws.Range("E32").Value = "  -1.45%  "
$ws.Range("D33").Value = "5.311"
$ws.Range("E33").Value = "  +1.04%  "
$ws.Range("D34").Value = "1.340"
$ws.Range("E34").Value = "  +1.57%  "
$ws.Range("D35").Value = "3.322"
$ws.Range("E35").Value = "  +0.44%  "
$ws.Range("D37").Value = "1.150"
$ws.Range("E37").Value = "  +0.72%  "
$ws.Range("D38").Value = "0.02071"
$ws.Range("E38").Value = "  -2.51%  "
$ws.Range("D39").Value = "7.675"
$ws.Range("E39").Value = "  -0.90%  "
$ws.Range("D40").Value = "0.5643"
$ws.Range("E40").Value = "  +0.98%  "
$ws.Range("E41").Value = "  +1.60%  "
$ws.Range("D42").Value = "9.836"
$ws.Range("E42").Value = "  -0.35%  "
$ws.Range("D43").Value = "0.07210"
$ws.Range("E43").Value = "  -0.07%  "
$ws.Range("D44").Value = "11.81"
$ws.Range("E44").Value = "  +1.93%  "
$ws.Range("D45").Value = "0.5296"
$ws.Range("E45").Value = "  +0.64%  "
$ws.Range("D46").Value = "2.137"
$ws.Range("E46").Value = "  +1.97%  "
$ws.Range("D47").Value = "1.126"
$ws.Range("E47").Value = "  -1.85%  "
$ws.Range("D48").Value = "1.831"
$ws.Range("E48").Value = "  +0.40%  "
$ws.Range("E49").Value = "  +0.48%  "
$ws.Range("D50").Value = "2.417"
$ws.Range("E50").Value = "  +4.11%  "
$ws.Range("E51").Value = "  +0.25%  "
